# TC_147 - add Alarm/Standby loading-detail summary columns (F:G) to the
# "Add Panels" sheet, mirroring the existing W:X columns, and fix up the
# row-3 height / column-G width / selection that changed along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# New header labels in F1:G1 - same look as the W7:X7 header cells
# (bold white-on-blue, bordered). Copy that formatting across so the
# engine reuses the existing style instead of inventing a new one.
$ws.Range("W7").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("F1").Value = "AlarmLoadingDetail"
$ws.Range("G1").Value = "StandbyLoadingDetail"

# New sub-labels in F2:G2 - bold, centered, bordered, no fill (same look
# as the existing C2 "Color Codes" cell).
$ws.Range("C2").Copy()
$ws.Range("F2:G2").PasteSpecial(-4122)
$ws.Range("F2").Value = "Alarm Current(A)"
$ws.Range("G2").Value = "Standby Current(A)"

# Row 3 no longer needs the taller 28.8pt height (back to sheet default).
$ws.Rows("3").RowHeight = 14.4

# Column G widened slightly to fit the new content.
$ws.Columns("G").ColumnWidth = 19.6640625

# Last selection when the file was saved.
$ws.Range("F4").Select()
